$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new value, and whether the value needs to be
# force-formatted as Text first (otherwise Excel auto-converts plain
# decimal-looking strings, e.g. "582.48", into a Number).
$updates = @(
    [PSCustomObject]@{ Cell = "D2"; Value = "63.783.54"; Force=$false }
    [PSCustomObject]@{ Cell = "E2"; Value = "  -0.77%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D3"; Value = "3.430.15"; Force=$false }
    [PSCustomObject]@{ Cell = "E3"; Value = "  -1.97%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D5"; Value = "582.48"; Force=$true }
    [PSCustomObject]@{ Cell = "E5"; Value = "  -1.01%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D6"; Value = "130.37"; Force=$true }
    [PSCustomObject]@{ Cell = "E6"; Value = "  -2.95%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E7"; Value = "  +0.05%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E8"; Value = "  -1.37%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D9"; Value = "7.58"; Force=$true }
    [PSCustomObject]@{ Cell = "E9"; Value = "  +4.05%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D10"; Value = "0.126"; Force=$true }
    [PSCustomObject]@{ Cell = "E10"; Value = "  +1.60%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E11"; Value = "  -0.85%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D12"; Value = "4.015.07"; Force=$false }
    [PSCustomObject]@{ Cell = "E12"; Value = "  -1.88%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E13"; Value = "  -0.28%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D14"; Value = "0.0000178"; Force=$true }
    [PSCustomObject]@{ Cell = "E14"; Value = "  -1.51%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D15"; Value = "3.430.05"; Force=$false }
    [PSCustomObject]@{ Cell = "E15"; Value = "  -1.88%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D16"; Value = "63.816.33"; Force=$false }
    [PSCustomObject]@{ Cell = "E16"; Value = "  -0.85%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D17"; Value = "25.02"; Force=$true }
    [PSCustomObject]@{ Cell = "E17"; Value = "  -2.44%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D18"; Value = "9.87"; Force=$true }
    [PSCustomObject]@{ Cell = "E18"; Value = "  +0.08%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E19"; Value = "  -1.20%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D20"; Value = "13.35"; Force=$true }
    [PSCustomObject]@{ Cell = "E20"; Value = "  -1.26%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E21"; Value = "  -2.06%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D22"; Value = "0.564"; Force=$true }
    [PSCustomObject]@{ Cell = "E22"; Value = "  -1.18%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D23"; Value = "3.567.65"; Force=$false }
    [PSCustomObject]@{ Cell = "E23"; Value = "  -1.89%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D24"; Value = "73.98"; Force=$true }
    [PSCustomObject]@{ Cell = "E24"; Value = "  -0.74%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E25"; Value = "  +0.12%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E26"; Value = "  -4.22%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D27"; Value = "0.998"; Force=$true }
    [PSCustomObject]@{ Cell = "E27"; Value = "  -0.16%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D28"; Value = "2.22"; Force=$true }
    [PSCustomObject]@{ Cell = "E28"; Value = "  -1.02%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D29"; Value = "7.06"; Force=$true }
    [PSCustomObject]@{ Cell = "E29"; Value = "  -4.08%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D30"; Value = "7.97"; Force=$true }
    [PSCustomObject]@{ Cell = "E30"; Value = "  -3.46%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E31"; Value = "  +1.99%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E32"; Value = "  -3.37%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D33"; Value = "3.458.67"; Force=$false }
    [PSCustomObject]@{ Cell = "E33"; Value = "  -1.76%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E34"; Value = "  -0.09%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D35"; Value = "22.95"; Force=$true }
    [PSCustomObject]@{ Cell = "E35"; Value = "  -2.12%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D36"; Value = "5.19"; Force=$true }
    [PSCustomObject]@{ Cell = "E36"; Value = "  +1.01%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E37"; Value = "  -1.45%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D38"; Value = "164.04"; Force=$true }
    [PSCustomObject]@{ Cell = "E38"; Value = "  -1.84%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E39"; Value = "  -2.05%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E40"; Value = "  -0.72%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D41"; Value = "0.788"; Force=$true }
    [PSCustomObject]@{ Cell = "E41"; Value = "  -2.67%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E42"; Value = "  +0.03%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D43"; Value = "41.47"; Force=$true }
    [PSCustomObject]@{ Cell = "E43"; Value = "  -1.07%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E44"; Value = "  -1.03%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E45"; Value = "  -2.18%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D46"; Value = "23.48"; Force=$true }
    [PSCustomObject]@{ Cell = "E46"; Value = "  -7.47%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E47"; Value = "  -4.62%  "; Force=$false }
    [PSCustomObject]@{ Cell = "E48"; Value = "  -0.18%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D49"; Value = "0.896"; Force=$true }
    [PSCustomObject]@{ Cell = "E49"; Value = "  +0.46%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D50"; Value = "2.288.61"; Force=$false }
    [PSCustomObject]@{ Cell = "E50"; Value = "  -7.24%  "; Force=$false }
    [PSCustomObject]@{ Cell = "D51"; Value = "0.0253"; Force=$true }
    [PSCustomObject]@{ Cell = "E51"; Value = "  -1.98%  "; Force=$false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Force) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
